$d = $word.ActiveDocument

# 1) Update the date in the header line and all verse texts in place
# via literal (non-wildcard) Find & Replace - every old verse text below
# is unique in the document, so each Execute() call touches exactly one run.
$replacements = @(
    ,@('********************************8月28日读经章节***************************', '********************************8月30日读经章节***************************')
    ,@('Chapter 1 of 1_Thessalonians', 'Chapter 3 of 1_Thessalonians')
    ,@('1.保罗，西拉，提摩太，写信给帖撒罗尼迦在父神和主耶稣基督里的教会。愿恩惠平安归与你们。', '1.我们既不能再忍，就愿意独自等在雅典。')
    ,@('2.我们为你们众人常常感谢神，祷告的时候提到你们。', '2.打发我们的兄弟在基督福音上作神执事的提摩太前去，（作神执事的有古卷作与神同工的）坚固你们，并在你们所信的道上劝慰你们。')
    ,@('3.在神我们的父面前，不住地记念你们因信心所作的工夫，因爱心所受的劳苦，因盼望我们主耶稣基督所存的忍耐。', '3.免得有人被诸般患难摇动。因为你们自己知道我们受患难原是命定的。')
    ,@('4.被神所爱的弟兄阿，我知道你们是蒙拣选的。', '4.我们在你们那里的时候，预先告诉你们，我们必受患难，以后果然应验了，你们也知道。')
    ,@('5.因为我们的福音传到你们那里，不独在乎言语，也在乎权能和圣灵，并充足的信心，正如你们知道我们在你们那里，为你们的缘故是怎样为人。', '5.为此，我既不能再忍，就打发人去，要晓得你们的信心如何，恐怕那诱惑人的到底诱惑了你们，叫我们的劳苦归于徒然。')
    ,@('6.并且你们在大难之中，蒙了圣灵所赐的喜乐，领受真道，就效法我们，也效法了主。', '6.但提摩太刚才从你们那里回来，将你们信心和爱心的好消息报给我们，又说你们常常记念我们，切切地想见我们，如同我们想见你们一样。')
    ,@('7.甚至你们作了马其顿和亚该亚，所有信主之人的榜样。', '7.所以弟兄们，我们在一切困苦患难之中，因着你们的信心就得了安慰。')
    ,@('8.因为主的道从你们那里已经传扬出来，你们向神的信心不但在马其顿和亚该亚，就是在各处，也都传开了。所以不用我们说什么话。', '8.你们若靠主站立得稳，我们就活了。')
    ,@('9.因为他们自己已经报明我们是怎样进到你们那里，你们是怎样离弃偶像归向神，要服事那又真又活的神，', '9.我们在神面前，因着你们甚是喜乐，为这一切喜乐，可用何等的感谢，为你们报答神呢？')
    ,@('10.等候他儿子从天降临，就是他从死里复活的，那位救我们脱离将来忿怒的耶稣。', '10.我们昼夜切切地祈求，要见你们的面，补满你们信心的不足。')
    ,@('Chapter 16 of Proverbs', '11.愿神我们的父，和我们的主耶稣，一直引领我们到你们那里去。')
    ,@('1.心中的谋算在乎人。舌头的应对，由于耶和华。', '12.又愿主叫你们彼此相爱的心，并爱众人的心，都能增长，充足，如同我们爱你们一样。')
    ,@('2.人一切所行的，在自己眼中看为清洁。惟有耶和华衡量人心。', '13.好使你们，当我们主耶稣同他众圣徒来的时候，在我们父神面前，心里坚固，成为圣洁，无可责备。')
    ,@('3.你所作的，要交托耶和华，你所谋的，就必成立。', 'Chapter 21 of Proverbs')
    ,@('4.耶和华所造的，各适其用。就是恶人，也为祸患的日子所造。', '1.王的心在耶和华手中，好像陇沟的水，随意流转。')
    ,@('5.凡心里骄傲的，为耶和华所憎恶。虽然连手，他必不免受罚。', '2.人所行的，在自己眼中都看为正，惟有耶和华衡量人心。')
    ,@('6.因怜悯诚实，罪孽得赎。敬畏耶和华的，远离恶事。', '3.行仁义公平，比献祭更蒙耶和华悦纳。')
    ,@('7.人所行的若蒙耶和华喜悦，耶和华也使他的仇敌与他和好。', '4.恶人发达，眼高心傲，这乃是罪。（发达原文作灯）')
    ,@('8.多有财利，行事不义，不如少有财利，行事公义。', '5.殷勤筹划的，足致丰裕。行事急躁的，都必缺乏。')
    ,@('9.人心筹算自己的道路。惟耶和华指引他的脚步。', '6.用诡诈之舌求财的，就是自己取死。所得之财，乃是吹来吹去的浮云。')
    ,@('10.王的嘴中有神语。审判之时，他的口，必不差错。', '7.恶人的强暴，必将自己扫除。因他们不肯按公平行事。')
    ,@('11.公道的天平和秤，都属耶和华。囊中一切法码，都为他所定。', '8.负罪之人的路，甚是弯曲。至于清洁的人，他所行的乃是正直。')
    ,@('12.作恶为王所憎恶。因国位是靠公义坚立。', '9.宁可住在房顶的角上，不在宽阔的房屋，与争吵的妇人同住。')
    ,@('13.公义的嘴，为王所喜悦。说正直话的，为王所喜爱。', '10.恶人的心，乐人受祸。他眼并不怜恤邻舍。')
    ,@('14.王的震怒，如杀人的使者。但智慧人能止息王怒。', '11.亵慢的人受刑罚，愚蒙的人就得智慧。智慧人受训诲，便得知识。')
    ,@('15.王的脸光，使人有生命。王的恩典，好像春云时雨。', '12.义人思想恶人的家，知道恶人倾倒，必至灭亡。')
    ,@('16.得智慧胜似得金子。选聪明强如选银子。', '13.塞耳不听穷人哀求的，他将来呼吁也不蒙应允。')
    ,@('17.正直人的道，是远离恶事。谨守己路的，是保全性命。', '14.暗中送的礼物，挽回怒气。怀中搋的贿赂，止息暴怒。')
    ,@('18.骄傲在败坏以先，狂心在跌倒之前。', '15.秉公行义，使义人喜乐，使作孽的人败坏。')
    ,@('19.心里谦卑与穷乏人来往，强如将掳物与骄傲人同分。', '16.迷离通达道路的，必住在阴魂的会中。')
    ,@('20.谨守训言的，必得好处。倚靠耶和华的，便为有福。', '17.爱宴乐的，必致穷乏。好酒爱膏油的，必不富足。')
    ,@('21.心中有智慧，必称为通达人。嘴中的甜言，加增人的学问。', '18.恶人作了义人的赎价。奸诈人代替正直人。')
    ,@('22.人有智慧就有生命的泉源。愚昧人必被愚昧惩治。', '19.宁可住在旷野，不与争吵使气的妇人同住。')
    ,@('23.智慧人的心，教训他的口，又使他的嘴，增长学问。', '20.智慧人家中积蓄宝物膏油。愚昧人随得来随吞下。')
    ,@('24.良言如同蜂房，使心觉甘甜，使骨得医治。', '21.追求公义仁慈的，就寻得生命，公义，和尊荣。')
    ,@('25.有一条路，人以为正，至终成为死亡之路。', '22.智慧人爬上勇士的城墙，倾覆他所倚靠的坚垒。')
    ,@('26.劳力人的胃口，使他劳力，因为他的口腹催逼他。', '23.谨守口与舌的，就保守自己免受灾难。')
    ,@('27.匪徒图谋奸恶，嘴上仿佛有烧焦的火。', '24.心骄气傲的人，名叫亵慢。他行事狂妄，都出于骄傲。')
    ,@('28.乖僻人播散分争。传舌的离间密友。', '25.懒惰人的心愿，将他杀害，因为他手不肯作工。')
    ,@('29.强暴人诱惑邻舍，领他走不善之道。', '26.有终日贪得无餍的，义人施舍而不吝惜。')
    ,@('30.眼目紧合的，图谋乖僻，嘴唇紧闭的，成就邪恶。', '27.恶人的祭物是可憎的，何况他存恶意来献呢？')
    ,@('31.白发是荣耀的冠冕。在公义的道上，必能得着。', '28.作假见证的必灭亡，惟有听真情而言的，其言长存。')
    ,@('32.不轻易发怒的，胜过勇士。治服己心的，强如取城。', '29.恶人脸无羞耻，正直人行事坚定。')
    ,@('33.签放在怀里。定事由耶和华。', '30.没有人能以智慧，聪明，谋略，敌挡耶和华。')
    ,@('Chapter 17 of Proverbs', '31.马是为打仗之日预备的。得胜乃在乎耶和华。')
    ,@('1.设筵满屋，大家相争，不如有块干饼，大家相安。', 'Chapter 22 of Proverbs')
    ,@('2.仆人办事聪明，必管辖贻羞之子，又在众子中，同分产业。', '1.美名胜过大财，恩宠强如金银。')
    ,@('3.鼎为炼银，炉为炼金。惟有耶和华熬炼人心。', '2.富户穷人，在世相遇，都为耶和华所造。')
    ,@('4.行恶的留心听奸诈之言。说谎的侧耳听邪恶之语。', '3.通达人见祸藏躲。愚蒙人前往受害。')
    ,@('5.戏笑穷人的，是辱没造他的主。幸灾乐祸的，必不免受罚。', '4.敬畏耶和华心存谦卑，就得富有，尊荣，生命，为赏赐。')
    ,@('6.子孙为老人的冠冕。父亲是儿女的荣耀。', '5.乖僻人的路上，有荆棘和网罗。保守自己生命的。必要远离。')
    ,@('7.愚顽人说美言本不相宜，何况君王说谎话呢？', '6.教养孩童，使他走当行的道，就是到老他也不偏离。')
    ,@('8.贿赂在馈送的人眼中，看为宝玉。随处运动，都得顺利。', '7.富户管辖穷人，欠债的是债主的仆人。')
    ,@('9.遮掩人过的，寻求人爱。屡次挑错的，离间密友。', '8.撒罪孽的，必收灾祸。他逞怒的杖，也必废掉。')
    ,@('10.一句责备话，深入聪明人的心，强如责打愚昧人一百下。', '9.眼目慈善的，就必蒙福。因他将食物分给穷人。')
    ,@('11.恶人只寻背叛，所以必有严厉的使者，奉差攻击他。', '10.赶出亵慢人，争端就消除，分争和羞辱，也必止息。')
    ,@('12.宁可遇见丢崽子的母熊，不可遇见正行愚妄的愚昧人。', '11.喜爱清心的人，因他嘴上的恩言，王必与他为友。')
    ,@('13.以恶报善的，祸患必不离他的家。', '12.耶和华的眼目，眷顾聪明人。却倾败奸诈人的言语。')
    ,@('14.分争的起头，如水放开。所以在争闹之先，必当止息争竞。', '13.懒惰人说，外头有狮子，我在街上，就必被杀。')
    ,@('15.定恶人为义的，定义人为恶的，这都为耶和华所憎恶。', '14.淫妇的口为深坑，耶和华所憎恶的，必陷在其中。')
    ,@('16.愚昧人既无聪明，为何手拿价银买智慧呢？', '15.愚蒙迷住孩童的心，用管教的杖可以远远赶除。')
    ,@('17.朋友乃时常亲爱。弟兄为患难而生。', '16.欺压贫穷为要利己的，并送礼与富户的，都必缺乏。')
    ,@('18.在邻舍面前击掌作保，乃是无知的人。', '17.你须侧耳听受智慧人的言语，留心领会我的知识，')
    ,@('19.喜爱争竞的，是喜爱过犯。高立家门的，乃自取败坏。', '18.你若心中存记，嘴上咬定，这便为美。')
    ,@('20.心存邪僻的，寻不着好处。舌弄是非的，陷在祸患中。', '19.我今日以此特特指教你，为要使你倚靠耶和华。')
    ,@('21.生愚昧子的，必自愁苦。愚顽人的父，毫无喜乐。', '20.谋略和知识的美事，我岂没有写给你吗？')
    ,@('22.喜乐的心，乃是良药。忧伤的灵，使骨枯干。', '21.要使你知道真言的实理，你好将真言回覆那打发你来的人。')
    ,@('23.恶人暗中受贿赂，为要颠倒判断。', '22.贫穷人，你不可因他贫穷，就抢夺他的物。也不可在城门口欺压困苦人。')
    ,@('24.明哲人眼前有智慧。愚昧人眼望地极。', '23.因耶和华必为他辨屈。抢夺他的，耶和华必夺取那人的命。')
    ,@('25.愚昧子使父亲愁烦，使母亲忧苦。', '24.好生气的人，不可与他结交。暴怒的人，不可与他来往。')
    ,@('26.刑罚义人为不善。责打君子为不义。', '25.恐怕你效法他的行为，自己就陷在网罗里。')
    ,@('27.寡少言语的有知识。性情温良的有聪明。', '26.不要与人击掌，不要为欠债的作保。')
    ,@('28.愚昧人若静默不言，也可算为智慧。闭口不说，也可算为聪明。', '27.你若没有什么偿还，何必使人夺去你睡卧的床呢？')
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $null = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# 2) Append the two new verses (22:28-29) at the very end of the reading,
#    after the existing trailing line break, each followed by its own
#    manual line break (Chr(11)), matching the document's <w:br/> style.
$lineBreak = [char]11

# Add a fresh trailing line break after the current last character,
# then insert the new verse text (with an internal line break) just
# before that new trailing break, so the final run ends on <w:br/> again.
$end = $d.Content.End
$tail = $d.Range($end - 1, $end - 1)
$tail.InsertAfter($lineBreak)

$end2 = $d.Content.End
$insertPoint = $d.Range($end2 - 2, $end2 - 2)
$insertPoint.InsertAfter('28.你先祖所立的地界，你不可挪移。' + $lineBreak + '29.你看见办事殷勤的人吗？他必站在君王面前，必不站在下贱人面前。')

